$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.658.90'
$ws.Range('E2').Value = '  -6.99%  '
$ws.Range('D3').Value = '1.696.01'
$ws.Range('E3').Value = '  -5.70%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '220.03'
$ws.Range('E5').Value = '  -4.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5129'
$ws.Range('E6').Value = '  -12.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.006'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2644'
$ws.Range('E8').Value = '  -4.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '22.17'
$ws.Range('E9').Value = '  -4.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06297'
$ws.Range('E10').Value = '  -7.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07340'
$ws.Range('E11').Value = '  -2.56%  '
$ws.Range('D12').Value = '1.704.00'
$ws.Range('E12').Value = '  -5.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.525'
$ws.Range('E13').Value = '  -5.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5795'
$ws.Range('E14').Value = '  -6.44%  '
$ws.Range('D15').Value = '1.927.86'
$ws.Range('E15').Value = '  -5.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008459'
$ws.Range('E16').Value = '  -7.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.52'
$ws.Range('D18').Value = '26.692.43'
$ws.Range('E18').Value = '  -6.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.005'
$ws.Range('E19').Value = '  -8.62%  '
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.00'
$ws.Range('E21').Value = '  -4.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '186.82'
$ws.Range('E22').Value = '  -11.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.254'
$ws.Range('E23').Value = '  -8.35%  '
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.78'
$ws.Range('E25').Value = '  -5.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.519'
$ws.Range('E26').Value = '  -5.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1158'
$ws.Range('E27').Value = '  -8.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.75'
$ws.Range('E28').Value = '  -4.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.354'
$ws.Range('E29').Value = '  -4.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05655'
$ws.Range('E30').Value = '  -7.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.340'
$ws.Range('E31').Value = '  -5.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.511'
$ws.Range('E32').Value = '  -7.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.494'
$ws.Range('E33').Value = '  -8.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.646'
$ws.Range('E34').Value = '  -5.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.020'
$ws.Range('E35').Value = '  -3.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6009'
$ws.Range('E36').Value = '  -6.55%  '
$ws.Range('E37').Value = '  -5.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.700'
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.102.02'
$ws.Range('E39').Value = '  -3.91%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01613'
$ws.Range('E40').Value = '  -4.96%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8600'
$ws.Range('E41').Value = '  -3.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.846'
$ws.Range('E42').Value = '  -10.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.004'
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.63'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('D45').Value = '1.854.18'
$ws.Range('E45').Value = '  -5.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000113'
$ws.Range('E46').Value = '  +1.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.73'
$ws.Range('E47').Value = '  -5.83%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.004'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.127'
$ws.Range('E49').Value = '  -3.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05243'
$ws.Range('E50').Value = '  -4.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4321'
$ws.Range('E51').Value = '  -3.55%  '
